# Duplicate Sheet1 (right after itself) so the new sheet inherits the same
# page setup / phonetic-guide / pageLayout view settings, then rename it.
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Copy($null, $ws1) | Out-Null
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Sheet2"

# Clear the copied data and put the single new string in A1, then advance
# the selection to A2 (as typing the value and pressing Enter would).
$ws2.Cells.Clear() | Out-Null
$ws2.Range("A1").Value = "bobfrompage2"
$ws2.Range("A2").Select() | Out-Null

# Add the cross-sheet formula on Sheet1!B2, re-activate Sheet1, and advance
# the selection to C2 (as typing the formula and pressing Tab would).
$ws1.Range("B2").Formula = "=Sheet2!A1"
$ws1.Activate() | Out-Null
$ws1.Range("C2").Select() | Out-Null
